# Sprint 2 burndown chart update
# - Row 10 ("Models") is renamed to "Dungeon Level Design" and its total
#   effort is bumped from 20 to 80, with days 6-10 (H:L) filled in.
# - Row 9 ("Animations") gets explicit 0s recorded for days 6-10 (H:L).
# - The dependent SUM/IFERROR formulas in rows 38/39 and the chart series
#   recalc automatically from these input changes.
# - Selection on "Current Iteration" moves from W12 to Z12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Current Iteration")

# Row 9 - "Animations": record explicit zero effort for days 6-10 (H9:L9)
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0

# Row 10 - rename "Models" to "Dungeon Level Design", bump estimate, log days 6-10
$ws.Range("A10").Value = "Dungeon Level Design"
$ws.Range("B10").Value = 80
$ws.Range("H10").Value = 13
$ws.Range("I10").Value = 6
$ws.Range("J10").Value = 3
$ws.Range("K10").Value = 4.5
$ws.Range("L10").Value = 5

# Update the active selection to Z12 (was W12)
$ws.Range("Z12").Select()
